$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing every existing
# row down by one. Excel's native Insert shifts cell values/styles but
# (matching the source data's generation pipeline) leaves the worksheet's
# <hyperlinks> ref->rId map untouched, so F2..F152 keep pointing at their
# original relationship ids after the shift.
$ws.Rows("2:2").Insert()

# Row 2 (freshly inserted) should carry the same Description/Product
# Code/Basic Price/Circular Date/Circular Link as the row now sitting at 3
# (which is what used to be row 2). Copy it down first so every column
# picks up the correct text, number formatting and style.
$ws.Range("A3:F3").Copy($ws.Range("A2:F2"))

# Only the Date column (A2) actually differs from row 3: it is the new
# day being published. Typing a dd-mm-yyyy-shaped string straight into a
# General-formatted cell gets auto-parsed into a date serial, so stage it
# in a scratch cell that is pre-formatted as Text, then copy only the
# value across (PasteSpecial values) so A2 keeps its original style/type.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "05-01-2026"
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)
$scratch.Clear()

# The newly inserted row at the bottom (153) is a duplicate of the old
# last row (152) - already correct after the Insert/shift above - but it
# still needs its own hyperlink relationship added (rId152), since that
# one genuinely is new.
$ws.Hyperlinks.Add($ws.Range("F153"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")
